{"js": "// Remove stray <w:lang w:val=\"en-US\"/> marks from a couple of paragraphs/runs\n// near the end of the document, and turn the (until now empty) paragraph that\n// follows \"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\" into a small bold Arial note\n// (\"\u041c\u0410\u041902146   4 \u043a\u043e\u0440\u043e\u0431\u043a\u0438\"), per the commit \"\u0438\u0448\u0434\u0432 19 - 030720\".\n//\n// Word's JS API does not expose run/paragraph language directly, so the\n// surgical way to drop a single <w:lang/> (while leaving every sibling run\n// untouched) is to replace the whole paragraph with an explicit OOXML\n// fragment that simply omits it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate \"\u0412\u0410\u0416\u041d\u041e: ...\" (paragraph-mark rPr + one run carry a now-unwanted\n// en-US language tag) and the two empty trailing paragraphs that sit right\n// after \"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\".\nlet vazhnoIndex = -1;\nlet limitIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (vazhnoIndex === -1 && t.indexOf(\"\u0412\u0410\u0416\u041d\u041e\") !== -1) {\n    vazhnoIndex = i;\n  }\n  if (t.indexOf(\"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\") !== -1) {\n    limitIndex = i;\n  }\n}\nif (vazhnoIndex === -1) {\n  throw new Error('Paragraph containing \"\u0412\u0410\u0416\u041d\u041e\" not found.');\n}\nif (limitIndex === -1) {\n  throw new Error('Paragraph containing \"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\" not found.');\n}\nconst firstEmptyIndex = limitIndex + 1;\nconst secondEmptyIndex = limitIndex + 2;\nif (secondEmptyIndex >= paragraphs.items.length) {\n  throw new Error(\"Expected two trailing empty paragraphs after the limit note.\");\n}\n\nconst vazhnoPara = paragraphs.items[vazhnoIndex];\nconst firstEmptyPara = paragraphs.items[firstEmptyIndex];\nconst secondEmptyPara = paragraphs.items[secondEmptyIndex];\n\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// 1) \"\u0412\u0410\u0416\u041d\u041e...\" paragraph: same runs/text as before, just without the two\n//    <w:lang w:val=\"en-US\"/> tags (one on the paragraph mark, one on the\n//    lone-space run right after \"\u041f\u043e\u0441\u0442\u0430\u0440\u0430\u044e\u0441\u044c\").\nconst vazhnoOoxml = wrapOoxml(\n  '<w:p>' +\n    '<w:pPr><w:rPr><w:color w:val=\"FF0000\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t>\u0412\u0410\u0416\u041d\u041e: \u0435\u0441\u043b\u0438 \u0432\u044b\u0439\u0442\u0438 \u0438\u0437 \u044d\u043a\u0440\u0430\u043d\u0430 \u0432\u044b\u043f\u043e\u043b\u043d\u0435\u043d\u0438\u044f \u043f\u043e\u0435\u0437\u0434\u043a\u0438 \u0438 \u0432\u0435\u0440\u043d\u0443\u0442\u044c\u0441\u044f \u043d\u0430\u0437\u0430\u0434, \u0442\u043e \u0441\u0447\u0435\u0442\u0447\u0438\u043a \u0434\u0438\u0441\u0442\u0430\u043d\u0446\u0438\u0438 \u043e\u0431\u043d\u0443\u043b\u044f\u0435\u0442\u0441\u044f</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t>. \u0414\u0430\u043d\u043d\u0430\u044f \u043f\u0440\u043e</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\">\u0431\u043b\u0435\u043c\u0430 \u043d\u0435 \u043c\u043e\u0436\u0435\u0442 \u0431\u044b\u0442\u044c </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t>\u0440\u0435\u0448\u0435\u043d\u0430</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\"> \u0442\u0430\u043a \u043a\u0430\u043a \u0441\u0447\u0435\u0442\u0447\u0438\u043a \u0432\u0440\u0435\u043c\u0435\u043d\u0438. \u041f\u043e\u0441\u0442\u0430\u0440\u0430\u044e\u0441\u044c</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n    '<w:t>\u0437\u0430\u0432\u0442\u0440\u0430 \u0437\u0430\u043a\u0440\u044b\u0442\u044c \u0434\u0430\u043d\u043d\u044b\u0435 \u043f\u0443\u043d\u043a\u0442\u044b</w:t></w:r>' +\n    '<w:r><w:rPr><w:color w:val=\"FF0000\"/></w:rPr><w:t>.</w:t></w:r>' +\n  '</w:p>'\n);\nvazhnoPara.insertOoxml(vazhnoOoxml, Word.InsertLocation.replace);\n\n// 2) First trailing empty paragraph becomes the new Arial/bold note.\nconst runRpr =\n  '<w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n  '<w:b/><w:bCs/><w:color w:val=\"000000\"/><w:sz w:val=\"9\"/><w:szCs w:val=\"9\"/>' +\n  '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>';\nconst newNoteOoxml = wrapOoxml(\n  '<w:p>' +\n    '<w:pPr><w:rPr>' + runRpr + '</w:rPr></w:pPr>' +\n    '<w:r><w:rPr>' + runRpr + '</w:rPr><w:t>\u041c\u0410\u041902146</w:t></w:r>' +\n    '<w:r><w:rPr>' + runRpr + '</w:rPr><w:t xml:space=\"preserve\">   4 </w:t></w:r>' +\n    '<w:r><w:rPr>' + runRpr + '</w:rPr><w:t>\u043a\u043e\u0440\u043e\u0431\u043a\u0438</w:t></w:r>' +\n  '</w:p>'\n);\nfirstEmptyPara.insertOoxml(newNoteOoxml, Word.InsertLocation.replace);\n\n// 3) Second trailing empty paragraph stays empty, just loses its <w:lang/>.\nconst secondEmptyOoxml = wrapOoxml(\n  '<w:p><w:pPr><w:rPr><w:color w:val=\"FF0000\"/></w:rPr></w:pPr></w:p>'\n);\nsecondEmptyPara.insertOoxml(secondEmptyOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Remove stray <w:lang w:val=\"en-US\"/> marks from a couple of paragraphs/runs\n# near the end of the document, and turn the (until now empty) paragraph that\n# follows \"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\" into a small bold Arial note\n# (\"\u041c\u0410\u041902146   4 \u043a\u043e\u0440\u043e\u0431\u043a\u0438\"), per the commit \"\u0438\u0448\u0434\u0432 19 - 030720\".\n#\n# The Word object model has no direct \"remove this run's language\" verb that\n# actually edits the OOXML here, so the surgical way to drop a single\n# <w:lang/> (while leaving every sibling run untouched) is to replace the\n# whole paragraph's Range contents via InsertXML with an explicit WordOpenXML\n# fragment that simply omits it.\n\n$d = $word.ActiveDocument\n\nfunction Wrap-Ooxml([string]$bodyXml) {\n    return @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>$bodyXml</w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n}\n\n# Locate the \"\u0412\u0410\u0416\u041d\u041e: ...\" paragraph and the \"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d\n# \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\" paragraph that precedes the two trailing empty paragraphs.\n$vazhnoIndex = -1\n$limitIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($vazhnoIndex -eq -1 -and $t.Contains(\"\u0412\u0410\u0416\u041d\u041e\")) {\n        $vazhnoIndex = $i\n    }\n    if ($t.Contains(\"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\")) {\n        $limitIndex = $i\n    }\n}\nif ($vazhnoIndex -eq -1) {\n    throw 'Paragraph containing \"\u0412\u0410\u0416\u041d\u041e\" not found.'\n}\nif ($limitIndex -eq -1) {\n    throw 'Paragraph containing \"\u041e\u0433\u0440\u0430\u043d\u0438\u0447\u0438\u0442\u044c \u043a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043c\u0430\u0448\u0438\u043d \u0440\u0430\u0434\u0438\u0443\u0441\u043e\u043c\" not found.'\n}\n$firstEmptyIndex = $limitIndex + 1\n$secondEmptyIndex = $limitIndex + 2\nif ($secondEmptyIndex -gt $d.Paragraphs.Count) {\n    throw \"Expected two trailing empty paragraphs after the limit note.\"\n}\n\n# 1) \"\u0412\u0410\u0416\u041d\u041e...\" paragraph: same runs/text as before, just without the two\n#    <w:lang w:val=\"en-US\"/> tags (one on the paragraph mark, one on the\n#    lone-space run right after \"\u041f\u043e\u0441\u0442\u0430\u0440\u0430\u044e\u0441\u044c\").\n$vazhnoBody = (\n    '<w:p>' +\n        '<w:pPr><w:rPr><w:color w:val=\"FF0000\"/></w:rPr></w:pPr>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t>\u0412\u0410\u0416\u041d\u041e: \u0435\u0441\u043b\u0438 \u0432\u044b\u0439\u0442\u0438 \u0438\u0437 \u044d\u043a\u0440\u0430\u043d\u0430 \u0432\u044b\u043f\u043e\u043b\u043d\u0435\u043d\u0438\u044f \u043f\u043e\u0435\u0437\u0434\u043a\u0438 \u0438 \u0432\u0435\u0440\u043d\u0443\u0442\u044c\u0441\u044f \u043d\u0430\u0437\u0430\u0434, \u0442\u043e \u0441\u0447\u0435\u0442\u0447\u0438\u043a \u0434\u0438\u0441\u0442\u0430\u043d\u0446\u0438\u0438 \u043e\u0431\u043d\u0443\u043b\u044f\u0435\u0442\u0441\u044f</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t>. \u0414\u0430\u043d\u043d\u0430\u044f \u043f\u0440\u043e</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t xml:space=\"preserve\">\u0431\u043b\u0435\u043c\u0430 \u043d\u0435 \u043c\u043e\u0436\u0435\u0442 \u0431\u044b\u0442\u044c </w:t></w:r>' +\n        '<w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t>\u0440\u0435\u0448\u0435\u043d\u0430</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t xml:space=\"preserve\"> \u0442\u0430\u043a \u043a\u0430\u043a \u0441\u0447\u0435\u0442\u0447\u0438\u043a \u0432\u0440\u0435\u043c\u0435\u043d\u0438. \u041f\u043e\u0441\u0442\u0430\u0440\u0430\u044e\u0441\u044c</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/><w:highlight w:val=\"yellow\"/></w:rPr>' +\n        '<w:t>\u0437\u0430\u0432\u0442\u0440\u0430 \u0437\u0430\u043a\u0440\u044b\u0442\u044c \u0434\u0430\u043d\u043d\u044b\u0435 \u043f\u0443\u043d\u043a\u0442\u044b</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"FF0000\"/></w:rPr><w:t>.</w:t></w:r>' +\n    '</w:p>'\n)\n$d.Paragraphs.Item($vazhnoIndex).Range.InsertXML((Wrap-Ooxml $vazhnoBody))\n\n# 2) First trailing empty paragraph becomes the new Arial/bold note.\n$runRpr = (\n    '<w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n    '<w:b/><w:bCs/><w:color w:val=\"000000\"/><w:sz w:val=\"9\"/><w:szCs w:val=\"9\"/>' +\n    '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/>'\n)\n$noteBody = (\n    '<w:p>' +\n        \"<w:pPr><w:rPr>$runRpr</w:rPr></w:pPr>\" +\n        \"<w:r><w:rPr>$runRpr</w:rPr><w:t>\u041c\u0410\u041902146</w:t></w:r>\" +\n        \"<w:r><w:rPr>$runRpr</w:rPr><w:t xml:space=`\"preserve`\">   4 </w:t></w:r>\" +\n        \"<w:r><w:rPr>$runRpr</w:rPr><w:t>\u043a\u043e\u0440\u043e\u0431\u043a\u0438</w:t></w:r>\" +\n    '</w:p>'\n)\n$d.Paragraphs.Item($firstEmptyIndex).Range.InsertXML((Wrap-Ooxml $noteBody))\n\n# 3) Second trailing empty paragraph stays empty, just loses its <w:lang/>.\n$emptyBody = '<w:p><w:pPr><w:rPr><w:color w:val=\"FF0000\"/></w:rPr></w:pPr></w:p>'\n$d.Paragraphs.Item($secondEmptyIndex).Range.InsertXML((Wrap-Ooxml $emptyBody))\n"}
